# Applies scheduled market-price refresh updates to the Leve profit sheets.
# Mirrors the upstream diff: per-row currentAveragePrice* / LevePrice* / LeveProfit*
# columns (H:N) are overwritten with freshly sampled market values.
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 3821
$ws.Range("I2").Value = 1000
$ws.Range("J2").Value = 5231.5
$ws.Range("K2").Value = 1000
$ws.Range("L2").Value = 5231.5
$ws.Range("M2").Value = -887
$ws.Range("N2").Value = -5457.5
$ws.Range("H29").Value = 91.5
$ws.Range("I29").Value = 91.5
$ws.Range("K29").Value = 274.5
$ws.Range("M29").Value = 6.5
$ws.Range("H32").Value = 1063.2858
$ws.Range("J32").Value = 1224.375
$ws.Range("L32").Value = 1224.375
$ws.Range("N32").Value = -1876.375
$ws.Range("H33").Value = 239
$ws.Range("I33").Value = 120.84615
$ws.Range("K33").Value = 120.84615
$ws.Range("M33").Value = 108.15385
$ws.Range("H40").Value = 1863.3334
$ws.Range("J40").Value = 1863.3334
$ws.Range("L40").Value = 1863.3334
$ws.Range("N40").Value = -2213.3334
$ws.Range("H135").Value = 468.33334
$ws.Range("I135").Value = 397
$ws.Range("K135").Value = 3573
$ws.Range("M135").Value = -1038
$ws.Range("H138").Value = 2692
$ws.Range("J138").Value = 2999.8235
$ws.Range("L138").Value = 8999.4705
$ws.Range("N138").Value = -19279.4705

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 34998.75
$ws.Range("J44").Value = 34998.75
$ws.Range("L44").Value = 34998.75
$ws.Range("N44").Value = -35974.75
$ws.Range("H62").Value = 44660
$ws.Range("J62").Value = 44660
$ws.Range("L62").Value = 44660
$ws.Range("N62").Value = -45908
$ws.Range("H65").Value = 44660
$ws.Range("J65").Value = 44660
$ws.Range("L65").Value = 133980
$ws.Range("N65").Value = -140220
$ws.Range("H132").Value = 2642.6667
$ws.Range("I132").Value = 2642.6667
$ws.Range("K132").Value = 7928.000100000001
$ws.Range("M132").Value = -5398.000100000001

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 16066.6
$ws.Range("J76").Value = 16066.6
$ws.Range("L76").Value = 16066.6
$ws.Range("N76").Value = -16696.6
$ws.Range("H79").Value = 16066.6
$ws.Range("J79").Value = 16066.6
$ws.Range("L79").Value = 16066.6
$ws.Range("N79").Value = -18250.6
$ws.Range("H80").Value = 518.63635
$ws.Range("J80").Value = 595.1429
$ws.Range("L80").Value = 595.1429
$ws.Range("N80").Value = -2591.1429
$ws.Range("H83").Value = 518.63635
$ws.Range("J83").Value = 595.1429
$ws.Range("L83").Value = 2975.7145
$ws.Range("N83").Value = -12959.7145
$ws.Range("H86").Value = 4099.2383
$ws.Range("I86").Value = 3627.7856
$ws.Range("J86").Value = 5042.143
$ws.Range("K86").Value = 3627.7856
$ws.Range("L86").Value = 5042.143
$ws.Range("M86").Value = -2504.7856
$ws.Range("N86").Value = -7288.143
$ws.Range("H89").Value = 4099.2383
$ws.Range("I89").Value = 3627.7856
$ws.Range("J89").Value = 5042.143
$ws.Range("K89").Value = 18138.928
$ws.Range("L89").Value = 25210.715
$ws.Range("M89").Value = -12522.928
$ws.Range("N89").Value = -36442.715
$ws.Range("H94").Value = 1254.1578
$ws.Range("I94").Value = 1301.6111
$ws.Range("K94").Value = 1301.6111
$ws.Range("M94").Value = -850.6111000000001

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 168.33333
$ws.Range("I10").Value = 200.4
$ws.Range("J10").Value = 8
$ws.Range("K10").Value = 200.4
$ws.Range("L10").Value = 8
$ws.Range("M10").Value = -61.40000000000001
$ws.Range("N10").Value = -286
$ws.Range("H60").Value = 22265.5
$ws.Range("I60").Value = 10093
$ws.Range("J60").Value = 24700
$ws.Range("K60").Value = 10093
$ws.Range("L60").Value = 24700
$ws.Range("M60").Value = -9582
$ws.Range("N60").Value = -25722
$ws.Range("H122").Value = 1383.1818
$ws.Range("I122").Value = 1271.5
$ws.Range("K122").Value = 3814.5
$ws.Range("M122").Value = -1364.5
$ws.Range("H132").Value = 1749
$ws.Range("I132").Value = 1749
$ws.Range("K132").Value = 5247
$ws.Range("M132").Value = -2717
$ws.Range("H134").Value = 3210.9
$ws.Range("I134").Value = 2513.75
$ws.Range("J134").Value = 5999.5
$ws.Range("K134").Value = 7541.25
$ws.Range("L134").Value = 17998.5
$ws.Range("M134").Value = -5006.25
$ws.Range("N134").Value = -23068.5

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 86.1
$ws.Range("J12").Value = 92.625
$ws.Range("L12").Value = 277.875
$ws.Range("N12").Value = -623.875
$ws.Range("H55").Value = 2895.7144
$ws.Range("J55").Value = 3041.5386
$ws.Range("L55").Value = 9124.6158
$ws.Range("N55").Value = -9478.6158
$ws.Range("H131").Value = 473.16666
$ws.Range("I131").Value = 473.16666
$ws.Range("K131").Value = 1419.49998
$ws.Range("M131").Value = 3620.50002
$ws.Range("H132").Value = 2300.5
$ws.Range("I132").Value = 2300.5
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 20704.5
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -18174.5
$ws.Range("N132").ClearContents()

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1539
$ws.Range("I16").Value = 1539
$ws.Range("K16").Value = 1539
$ws.Range("M16").Value = -1369
$ws.Range("H40").Value = 1958.25
$ws.Range("I40").Value = 1889.3334
$ws.Range("J40").Value = 2165
$ws.Range("K40").Value = 1889.3334
$ws.Range("L40").Value = 2165
$ws.Range("M40").Value = -1753.3334
$ws.Range("N40").Value = -2437

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 25214.5
$ws.Range("J46").Value = 25214.5
$ws.Range("L46").Value = 25214.5
$ws.Range("N46").Value = -25676.5
$ws.Range("H81").Value = 2244
$ws.Range("I81").Value = 2071.2222
$ws.Range("K81").Value = 4142.4444
$ws.Range("M81").Value = -3081.4444
$ws.Range("H84").Value = 2244
$ws.Range("I84").Value = 2071.2222
$ws.Range("K84").Value = 20712.222
$ws.Range("M84").Value = -15408.222
$ws.Range("H98").Value = 13225
$ws.Range("J98").Value = 13225
$ws.Range("L98").Value = 13225
$ws.Range("N98").Value = -19215
$ws.Range("H126").Value = 1037.25
$ws.Range("I126").Value = 763.5455
$ws.Range("K126").Value = 2290.6365
$ws.Range("M126").Value = 179.3635000000004
$ws.Range("H134").Value = 25214.5
$ws.Range("J134").Value = 25214.5
$ws.Range("L134").Value = 75643.5
$ws.Range("N134").Value = -80713.5

